$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '92.235.78'
$ws.Range("E2").Value = '  +6.35%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.284.27'
$ws.Range("E3").Value = '  +1.23%  '

# Row 4
$ws.Range("E4").Value = '  +0.10%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.12'
$ws.Range("E5").Value = '  +1.95%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.04'
$ws.Range("E6").Value = '  +1.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.416'
$ws.Range("E7").Value = '  +15.58%  '

# Row 8
$ws.Range("E8").Value = '  +5.51%  '

# Row 9
$ws.Range("E9").Value = '  -0.08%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.282.71'
$ws.Range("E10").Value = '  +1.31%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.585'
$ws.Range("E11").Value = '  +3.09%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000265'
$ws.Range("E12").Value = '  +5.35%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.180'
$ws.Range("E13").Value = '  +0.93%  '

# Row 14
$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.909.02'
$ws.Range("E14").Value = '  +1.86%  '

# Row 15
$ws.Range("B15").Value = 'Avalanche'
$ws.Range("C15").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.28'
$ws.Range("E15").Value = '  +1.97%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '92.148.66'
$ws.Range("E16").Value = '  +6.56%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.36'
$ws.Range("E17").Value = '  +1.53%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.311.33'
$ws.Range("E18").Value = '  +2.29%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.27'
$ws.Range("E19").Value = '  +6.54%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.07'
$ws.Range("E20").Value = '  +1.35%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '438.92'
$ws.Range("E21").Value = '  +2.35%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.87'
$ws.Range("E22").Value = '  +1.05%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.28'
$ws.Range("E23").Value = '  -0.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000187'
$ws.Range("E24").Value = '  +48.52%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.36'
$ws.Range("E25").Value = '  +6.21%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.26'
$ws.Range("E26").Value = '  -1.29%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.528.88'
$ws.Range("E27").Value = '  +3.37%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '76.31'
$ws.Range("E28").Value = '  +1.67%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  +0.01%  '

# Row 30
$ws.Range("E30").Value = '  +4.03%  '

# Row 31
$ws.Range("E31").Value = '  +0.32%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.84'
$ws.Range("E32").Value = '  +1.78%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '552.77'
$ws.Range("E33").Value = '  +2.76%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.16'
$ws.Range("E34").Value = '  +3.56%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.64'
$ws.Range("E35").Value = '  +26.12%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.92'
$ws.Range("E36").Value = '  -0.13%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.30'
$ws.Range("E37").Value = '  -6.57%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '22.74'
$ws.Range("E38").Value = '  +2.62%  '

# Row 39
$ws.Range("E39").Value = '  -2.64%  '

# Row 40
$ws.Range("E40").Value = '  +3.96%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.01%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.394'
$ws.Range("E42").Value = '  +2.01%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.98'
$ws.Range("E43").Value = '  +0.25%  '

# Row 44
$ws.Range("E44").Value = '  +0.18%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.48'
$ws.Range("E45").Value = '  -2.28%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '181.07'
$ws.Range("E46").Value = '  +2.12%  '

# Row 47
$ws.Range("B47").Value = 'OKB'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '43.71'
$ws.Range("E47").Value = '  -1.26%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.130'
$ws.Range("E48").Value = '  +6.89%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.28'
$ws.Range("E49").Value = '  +0.33%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.631'
$ws.Range("E50").Value = '  +2.41%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.17'
$ws.Range("E51").Value = '  +0.13%  '
